$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2506.2
$ws.Range("J32").Value = 2506.2
$ws.Range("L32").Value = 2506.2
$ws.Range("N32").Value = -3158.2
$ws.Range("H54").Value = 5000
$ws.Range("I54").Value = 5000
$ws.Range("K54").Value = 5000
$ws.Range("M54").Value = -4514
$ws.Range("H62").Value = 100037896
$ws.Range("I62").Value = 200003680
$ws.Range("K62").Value = 200003680
$ws.Range("M62").Value = -200003056
$ws.Range("H65").Value = 100037896
$ws.Range("I65").Value = 200003680
$ws.Range("K65").Value = 1000018400
$ws.Range("M65").Value = -1000015280
$ws.Range("H92").Value = 998.3333
$ws.Range("I92").Value = 585.25
$ws.Range("J92").Value = 1824.5
$ws.Range("K92").Value = 585.25
$ws.Range("L92").Value = 1824.5
$ws.Range("M92").Value = 662.75
$ws.Range("N92").Value = -4320.5
$ws.Range("H100").Value = 3018.7778
$ws.Range("I100").Value = 1866.6666
$ws.Range("J100").Value = 5323
$ws.Range("K100").Value = 1866.6666
$ws.Range("L100").Value = 5323
$ws.Range("M100").Value = -1325.6666
$ws.Range("N100").Value = -6405
$ws.Range("H127").Value = 2523.5557
$ws.Range("I127").Value = 970.7143
$ws.Range("K127").Value = 2912.1429
$ws.Range("M127").Value = 2047.8571
$ws.Range("H138").Value = 5941.59
$ws.Range("J138").Value = 11719.889
$ws.Range("L138").Value = 35159.667
$ws.Range("N138").Value = -45439.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1931.5
$ws.Range("I16").Value = 897.25
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 897.25
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -610.25
$ws.Range("N16").Value = -4574
$ws.Range("H32").Value = 2781029.8
$ws.Range("I32").Value = 2844076.5
$ws.Range("K32").Value = 2844076.5
$ws.Range("M32").Value = -2843789.5
$ws.Range("H74").Value = 25011.615
$ws.Range("I74").Value = 30840.648
$ws.Range("K74").Value = 30840.648
$ws.Range("M74").Value = -29966.648
$ws.Range("H77").Value = 25011.615
$ws.Range("I77").Value = 30840.648
$ws.Range("K77").Value = 154203.24
$ws.Range("M77").Value = -149835.24
$ws.Range("H97").Value = 4173125.2
$ws.Range("J97").Value = 11922274
$ws.Range("L97").Value = 11922274
$ws.Range("N97").Value = -11923266
$ws.Range("H122").Value = 16846.941
$ws.Range("I122").Value = 20783.334
$ws.Range("K122").Value = 62350.00199999999
$ws.Range("M122").Value = -59900.00199999999
$ws.Range("H132").Value = 4638.1665
$ws.Range("I132").Value = 4023.8823
$ws.Range("K132").Value = 12071.6469
$ws.Range("M132").Value = -9541.6469

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1481.6666
$ws.Range("I94").Value = 538.82355
$ws.Range("K94").Value = 538.82355
$ws.Range("M94").Value = -87.82354999999995
$ws.Range("H105").Value = 3009.7837
$ws.Range("I105").Value = 2582.7778
$ws.Range("J105").Value = 4162.7
$ws.Range("K105").Value = 2582.7778
$ws.Range("L105").Value = 4162.7
$ws.Range("M105").Value = -835.7777999999998
$ws.Range("N105").Value = -7656.7
$ws.Range("H113").Value = 5013
$ws.Range("I113").Value = 5013
$ws.Range("K113").Value = 5013
$ws.Range("M113").Value = -2843

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H31").Value = 6838.25
$ws.Range("I31").Value = 5130
$ws.Range("J31").Value = 7358.1523
$ws.Range("K31").Value = 5130
$ws.Range("L31").Value = 7358.1523
$ws.Range("M31").Value = -4835
$ws.Range("N31").Value = -7948.1523
$ws.Range("H34").Value = 6838.25
$ws.Range("I34").Value = 5130
$ws.Range("J34").Value = 7358.1523
$ws.Range("K34").Value = 5130
$ws.Range("L34").Value = 7358.1523
$ws.Range("M34").Value = -4928
$ws.Range("N34").Value = -7762.1523
$ws.Range("H58").Value = 11909446
$ws.Range("I58").Value = 22729746
$ws.Range("J58").Value = 7116
$ws.Range("K58").Value = 22729746
$ws.Range("L58").Value = 7116
$ws.Range("M58").Value = -22729543
$ws.Range("N58").Value = -7522
$ws.Range("H99").Value = 10061
$ws.Range("I99").Value = 11135.167
$ws.Range("J99").Value = 8449.75
$ws.Range("K99").Value = 11135.167
$ws.Range("L99").Value = 8449.75
$ws.Range("M99").Value = -9637.166999999999
$ws.Range("N99").Value = -11445.75
$ws.Range("H107").Value = 1903.125
$ws.Range("I107").Value = 1367.75
$ws.Range("K107").Value = 1367.75
$ws.Range("M107").Value = 552.25
$ws.Range("H126").Value = 10061
$ws.Range("I126").Value = 11135.167
$ws.Range("J126").Value = 8449.75
$ws.Range("K126").Value = 33405.501
$ws.Range("L126").Value = 25349.25
$ws.Range("M126").Value = -30935.501
$ws.Range("N126").Value = -30289.25
$ws.Range("H134").Value = 17863200
$ws.Range("I134").Value = 6189.8184
$ws.Range("K134").Value = 18569.4552
$ws.Range("M134").Value = -16034.4552
$ws.Range("H136").Value = 11909446
$ws.Range("I136").Value = 22729746
$ws.Range("J136").Value = 7116
$ws.Range("K136").Value = 68189238
$ws.Range("L136").Value = 21348
$ws.Range("M136").Value = -68186688
$ws.Range("N136").Value = -26448
$ws.Range("H141").Value = 55414.875
$ws.Range("J141").Value = 53885.332
$ws.Range("L141").Value = 53885.332
$ws.Range("N141").Value = -64245.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 15151726
$ws.Range("I14").Value = 15151726
$ws.Range("K14").Value = 45455178
$ws.Range("M14").Value = -45455005
$ws.Range("H80").Value = 250003000
$ws.Range("I80").Value = 250004000
$ws.Range("K80").Value = 750012000
$ws.Range("M80").Value = -750011064
$ws.Range("H83").Value = 250003000
$ws.Range("I83").Value = 250004000
$ws.Range("K83").Value = 2250036000
$ws.Range("M83").Value = -2250031320
$ws.Range("H87").Value = 12000
$ws.Range("J87").Value = 12000
$ws.Range("L87").Value = 36000
$ws.Range("N87").Value = -38496
$ws.Range("H90").Value = 12000
$ws.Range("J90").Value = 12000
$ws.Range("L90").Value = 108000
$ws.Range("N90").Value = -120480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 39647.332
$ws.Range("J95").Value = 39647.332
$ws.Range("L95").Value = 39647.332
$ws.Range("N95").Value = -45139.332
$ws.Range("H122").Value = 51144.863
$ws.Range("I122").Value = 95474
$ws.Range("J122").Value = 6815.727
$ws.Range("K122").Value = 286422
$ws.Range("L122").Value = 20447.181
$ws.Range("M122").Value = -283972
$ws.Range("N122").Value = -25347.181
$ws.Range("H126").Value = 3485.2856
$ws.Range("I126").Value = 3399
$ws.Range("J126").Value = 3499.6667
$ws.Range("K126").Value = 10197
$ws.Range("L126").Value = 10499.0001
$ws.Range("M126").Value = -7727
$ws.Range("N126").Value = -15439.0001
$ws.Range("H131").Value = 59520.5
$ws.Range("J131").Value = 59520.5
$ws.Range("L131").Value = 59520.5
$ws.Range("N131").Value = -69600.5
$ws.Range("H132").Value = 2403.1794
$ws.Range("I132").Value = 2249.9092
$ws.Range("J132").Value = 3246.1667
$ws.Range("K132").Value = 6749.7276
$ws.Range("L132").Value = 9738.500100000001
$ws.Range("M132").Value = -4219.7276
$ws.Range("N132").Value = -14798.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7078.579
$ws.Range("I7").Value = 5198.1665
$ws.Range("K7").Value = 5198.1665
$ws.Range("M7").Value = -5086.1665
$ws.Range("H82").Value = 2014.8823
$ws.Range("I82").Value = 1985.2222
$ws.Range("J82").Value = 2048.25
$ws.Range("K82").Value = 1985.2222
$ws.Range("L82").Value = 2048.25
$ws.Range("M82").Value = -1624.2222
$ws.Range("N82").Value = -2770.25
$ws.Range("H85").Value = 2014.8823
$ws.Range("I85").Value = 1985.2222
$ws.Range("J85").Value = 2048.25
$ws.Range("K85").Value = 1985.2222
$ws.Range("L85").Value = 2048.25
$ws.Range("M85").Value = -737.2221999999999
$ws.Range("N85").Value = -4544.25
$ws.Range("H100").Value = 3195.9583
$ws.Range("I100").Value = 2785.75
$ws.Range("J100").Value = 3401.0625
$ws.Range("K100").Value = 2785.75
$ws.Range("L100").Value = 3401.0625
$ws.Range("M100").Value = -2244.75
$ws.Range("N100").Value = -4483.0625
$ws.Range("H125").Value = 65251
$ws.Range("J125").Value = 65251
$ws.Range("L125").Value = 65251
$ws.Range("N125").Value = -75091
$ws.Range("H126").Value = 7078.579
$ws.Range("I126").Value = 5198.1665
$ws.Range("K126").Value = 15594.4995
$ws.Range("M126").Value = -13124.4995
$ws.Range("H132").Value = 22736980
$ws.Range("I132").Value = 45466000
$ws.Range("K132").Value = 136398000
$ws.Range("M132").Value = -136395470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3759.8667
$ws.Range("I126").Value = 3233.1667
$ws.Range("J126").Value = 4111
$ws.Range("K126").Value = 9699.500100000001
$ws.Range("L126").Value = 12333
$ws.Range("M126").Value = -7229.500100000001
$ws.Range("N126").Value = -17273
$ws.Range("H132").Value = 21758228
$ws.Range("I132").Value = 26324756
$ws.Range("K132").Value = 78974268
$ws.Range("M132").Value = -78971738
$ws.Range("H135").Value = 84999
$ws.Range("J135").Value = 84999
$ws.Range("L135").Value = 84999
$ws.Range("N135").Value = -95139
$ws.Range("H136").Value = 69004270
$ws.Range("I136").Value = 133334616
$ws.Range("K136").Value = 400003848
$ws.Range("M136").Value = -400001298
